$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.059.43"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").Value = "2.176.49"
$ws.Range("E3").Value = "  -2.72%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.49"
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.619"
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "67.20"
$ws.Range("E7").Value = "  -6.64%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.569"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.20"
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0928"
$ws.Range("E11").Value = "  -4.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "35.61"
$ws.Range("E12").Value = "  -16.70%  "
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.91"
$ws.Range("E14").Value = "  -1.26%  "
$ws.Range("D15").Value = "2.500.40"
$ws.Range("E15").Value = "  -2.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.857"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.22"
$ws.Range("E17").Value = "  -5.83%  "
$ws.Range("D18").Value = "2.179.24"
$ws.Range("E18").Value = "  -2.57%  "
$ws.Range("D19").Value = "41.001.49"
$ws.Range("E19").Value = "  -1.88%  "
$ws.Range("D20").Value = "0.0₃0941"
$ws.Range("E20").Value = "  -2.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.09"
$ws.Range("E21").Value = "  -2.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.63"
$ws.Range("E22").Value = "  -2.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.01"
$ws.Range("E23").Value = "  -2.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.06"
$ws.Range("E24").Value = "  -9.96%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.51"
$ws.Range("E25").Value = "  +11.77%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("B27").Value = "WEMIXToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.74"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.42"
$ws.Range("E28").Value = "  -3.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "168.62"
$ws.Range("E29").Value = "  -1.98%  "
$ws.Range("E30").Value = "  -6.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.25"
$ws.Range("E31").Value = "  -2.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.120"
$ws.Range("E32").Value = "  -2.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.64"
$ws.Range("E33").Value = "  +2.91%  "
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("E35").Value = "  -3.69%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.10"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.57"
$ws.Range("E37").Value = "  -3.17%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.53"
$ws.Range("E38").Value = "  -4.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0301"
$ws.Range("E39").Value = "  +7.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.17"
$ws.Range("E40").Value = "  -5.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.50"
$ws.Range("E41").Value = "  -9.73%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.50"
$ws.Range("E42").Value = "  -2.44%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "60.94"
$ws.Range("E43").Value = "  -11.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.77"
$ws.Range("E44").Value = "  -7.28%  "
$ws.Range("E45").Value = "  -10.85%  "
$ws.Range("E46").Value = "  -5.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.01"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0987"
$ws.Range("E48").Value = "  -3.56%  "
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.16"
$ws.Range("E50").Value = "  -3.65%  "
$ws.Range("E51").Value = "  -0.44%  "
